# Commit: "Changed and updated test cases for Bootstrap"
#
# - Bootstrap sheet (row 3-8): add "Same as expected" / "Pass" result
#   columns (G/H), edit the error-message text in F4, and drop the
#   leftover SUM formula that used to live in G8.
# - AddBatch sheet: loses the "active tab" marker (Bootstrap becomes the
#   active tab/sheet instead), shared-string indices shift as a natural
#   consequence of the sharedStrings table being de-duplicated on save.

$wb = $excel.ActiveWorkbook
$bootstrap = $wb.Worksheets.Item("Bootstrap")

# --- Row 4's F4 description text is edited first (row 5 dropped from the
#     list of rows affected by the password error) so the edited string
#     is registered in the shared-string table ahead of "Same as expected".
$bootstrap.Range("F4").Value = "Bootstrap success. However, 8 error messages are returned to admin. Row 3,4 will consist of mac-address error. Row 6,7 will consist of password error. Row 8,9 will consist of email error. Row 10 will consist of gender error. Rows 3,4,6,7,8,9,10 will not be displayed in the database. All other rows should be displayed correctly."

# --- Row 3: G3 keeps its existing style (vertical-center only), H3 gets
#     the same centered/wrapped style already used across the table body.
$bootstrap.Range("G3").Value = "Same as expected"

$bootstrap.Range("H3").HorizontalAlignment = -4108
$bootstrap.Range("H3").VerticalAlignment = -4108
$bootstrap.Range("H3").WrapText = $true
$bootstrap.Range("H3").Value = "Pass"

# --- Row 4: G4/H4 are brand new plain cells (no explicit style),
#     recording the same "Same as expected" / "Pass" outcome.
$bootstrap.Range("G4").Value = "Same as expected"
$bootstrap.Range("H4").Value = "Pass"

# --- Row 5: same pattern as row 4 (plain, unstyled new cells).
$bootstrap.Range("G5").Value = "Same as expected"
$bootstrap.Range("H5").Value = "Pass"

# --- Row 6: new cells use the centered/wrapped body style.
$bootstrap.Range("G6").HorizontalAlignment = -4108
$bootstrap.Range("G6").VerticalAlignment = -4108
$bootstrap.Range("G6").WrapText = $true
$bootstrap.Range("G6").Value = "Same as expected"

$bootstrap.Range("H6").HorizontalAlignment = -4108
$bootstrap.Range("H6").VerticalAlignment = -4108
$bootstrap.Range("H6").WrapText = $true
$bootstrap.Range("H6").Value = "Pass"

# --- Row 7: same centered/wrapped body style.
$bootstrap.Range("G7").HorizontalAlignment = -4108
$bootstrap.Range("G7").VerticalAlignment = -4108
$bootstrap.Range("G7").WrapText = $true
$bootstrap.Range("G7").Value = "Same as expected"

$bootstrap.Range("H7").HorizontalAlignment = -4108
$bootstrap.Range("H7").VerticalAlignment = -4108
$bootstrap.Range("H7").WrapText = $true
$bootstrap.Range("H7").Value = "Pass"

# --- Row 8: G8 used to hold a leftover "=SUM(33-5)" formula; replace it
#     with the same "Same as expected" text (centered/wrapped style),
#     and add H8 ("Pass") with the same style.
$bootstrap.Range("G8").HorizontalAlignment = -4108
$bootstrap.Range("G8").VerticalAlignment = -4108
$bootstrap.Range("G8").WrapText = $true
$bootstrap.Range("G8").Value = "Same as expected"

$bootstrap.Range("H8").HorizontalAlignment = -4108
$bootstrap.Range("H8").VerticalAlignment = -4108
$bootstrap.Range("H8").WrapText = $true
$bootstrap.Range("H8").Value = "Pass"

# --- Make Bootstrap the active/selected sheet and select G8, which also
#     clears the previous "tabSelected" flag on AddBatch and moves the
#     workbook's activeTab pointer.
$bootstrap.Select()
$bootstrap.Range("G8").Select()
